$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value2 = 3946.3635
$ws.Cells.Item(62, 9).Value2 = 3652.5
$ws.Cells.Item(62, 10).Value2 = 4114.2856
$ws.Cells.Item(62, 11).Value2 = 3652.5
$ws.Cells.Item(62, 12).Value2 = 4114.2856
$ws.Cells.Item(62, 13).Value2 = -3028.5
$ws.Cells.Item(62, 14).Value2 = -5362.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value2 = 3946.3635
$ws.Cells.Item(65, 9).Value2 = 3652.5
$ws.Cells.Item(65, 10).Value2 = 4114.2856
$ws.Cells.Item(65, 11).Value2 = 18262.5
$ws.Cells.Item(65, 12).Value2 = 20571.428
$ws.Cells.Item(65, 13).Value2 = -15142.5
$ws.Cells.Item(65, 14).Value2 = -26811.428

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(87, 8).Value2 = 34108.785
$ws.Cells.Item(87, 10).Value2 = 34108.785
$ws.Cells.Item(87, 12).Value2 = 34108.785
$ws.Cells.Item(87, 14).Value2 = -36604.785

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value2 = 1541.8572
$ws.Cells.Item(88, 9).Value2 = 1446.5
$ws.Cells.Item(88, 10).Value2 = 1580
$ws.Cells.Item(88, 11).Value2 = 1446.5
$ws.Cells.Item(88, 12).Value2 = 1580
$ws.Cells.Item(88, 13).Value2 = -1040.5
$ws.Cells.Item(88, 14).Value2 = -2392

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(90, 8).Value2 = 34108.785
$ws.Cells.Item(90, 10).Value2 = 34108.785
$ws.Cells.Item(90, 12).Value2 = 102326.355
$ws.Cells.Item(90, 14).Value2 = -114806.355

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value2 = 1541.8572
$ws.Cells.Item(91, 9).Value2 = 1446.5
$ws.Cells.Item(91, 10).Value2 = 1580
$ws.Cells.Item(91, 11).Value2 = 1446.5
$ws.Cells.Item(91, 12).Value2 = 1580
$ws.Cells.Item(91, 13).Value2 = -42.5
$ws.Cells.Item(91, 14).Value2 = -4388

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value2 = 1242.6207
$ws.Cells.Item(112, 10).Value2 = 1242.6207
$ws.Cells.Item(112, 12).Value2 = 3727.8621
$ws.Cells.Item(112, 14).Value2 = -5943.8621

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value2 = 1016.3148
$ws.Cells.Item(129, 9).Value2 = 482.33334
$ws.Cells.Item(129, 10).Value2 = 1047.7255
$ws.Cells.Item(129, 11).Value2 = 1447.00002
$ws.Cells.Item(129, 12).Value2 = 3143.1765
$ws.Cells.Item(129, 13).Value2 = 3552.99998
$ws.Cells.Item(129, 14).Value2 = -13143.1765

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value2 = 2301.6956
$ws.Cells.Item(137, 9).Value2 = 1143.5957
$ws.Cells.Item(137, 10).Value2 = 4775.8184
$ws.Cells.Item(137, 11).Value2 = 3430.7871
$ws.Cells.Item(137, 12).Value2 = 14327.4552
$ws.Cells.Item(137, 13).Value2 = -880.7871000000005
$ws.Cells.Item(137, 14).Value2 = -19427.4552

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value2 = 3378.2158
$ws.Cells.Item(138, 9).Value2 = 1774.64
$ws.Cells.Item(138, 10).Value2 = 4014.5557
$ws.Cells.Item(138, 11).Value2 = 5323.92
$ws.Cells.Item(138, 12).Value2 = 12043.6671
$ws.Cells.Item(138, 13).Value2 = -183.9200000000001
$ws.Cells.Item(138, 14).Value2 = -22323.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value2 = 6312.8887
$ws.Cells.Item(61, 9).Value2 = 4632.45
$ws.Cells.Item(61, 11).Value2 = 4632.45
$ws.Cells.Item(61, 13).Value2 = -4420.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value2 = 6832.963
$ws.Cells.Item(74, 9).Value2 = 3969.7646
$ws.Cells.Item(74, 10).Value2 = 11700.4
$ws.Cells.Item(74, 11).Value2 = 3969.7646
$ws.Cells.Item(74, 12).Value2 = 11700.4
$ws.Cells.Item(74, 13).Value2 = -3095.7646
$ws.Cells.Item(74, 14).Value2 = -13448.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value2 = 6832.963
$ws.Cells.Item(77, 9).Value2 = 3969.7646
$ws.Cells.Item(77, 10).Value2 = 11700.4
$ws.Cells.Item(77, 11).Value2 = 19848.823
$ws.Cells.Item(77, 12).Value2 = 58502
$ws.Cells.Item(77, 13).Value2 = -15480.823
$ws.Cells.Item(77, 14).Value2 = -67238

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value2 = 5102.814
$ws.Cells.Item(132, 9).Value2 = 1434
$ws.Cells.Item(132, 10).Value2 = 8293.087
$ws.Cells.Item(132, 11).Value2 = 4302
$ws.Cells.Item(132, 12).Value2 = 24879.261
$ws.Cells.Item(132, 13).Value2 = -1772
$ws.Cells.Item(132, 14).Value2 = -29939.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value2 = 6312.8887
$ws.Cells.Item(136, 9).Value2 = 4632.45
$ws.Cells.Item(136, 11).Value2 = 13897.35
$ws.Cells.Item(136, 13).Value2 = -11347.35

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value2 = 5647.2666
$ws.Cells.Item(105, 9).Value2 = 5577.6665
$ws.Cells.Item(105, 10).Value2 = 5751.6665
$ws.Cells.Item(105, 11).Value2 = 5577.6665
$ws.Cells.Item(105, 12).Value2 = 5751.6665
$ws.Cells.Item(105, 13).Value2 = -3830.6665
$ws.Cells.Item(105, 14).Value2 = -9245.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value2 = 2004.875
$ws.Cells.Item(107, 10).Value2 = 1904.6666
$ws.Cells.Item(107, 12).Value2 = 1904.6666
$ws.Cells.Item(107, 14).Value2 = -5744.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value2 = 5259.8125
$ws.Cells.Item(134, 9).Value2 = 4741.9585
$ws.Cells.Item(134, 10).Value2 = 6813.375
$ws.Cells.Item(134, 11).Value2 = 14225.8755
$ws.Cells.Item(134, 12).Value2 = 20440.125
$ws.Cells.Item(134, 13).Value2 = -11690.8755
$ws.Cells.Item(134, 14).Value2 = -25510.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value2 = 2930
$ws.Cells.Item(99, 9).Value2 = 2930
$ws.Cells.Item(99, 10).Value2 = 0
$ws.Cells.Item(99, 11).Value2 = 2930
$ws.Cells.Item(99, 12).Value2 = 0
$ws.Cells.Item(99, 13).Value2 = -1432
$ws.Cells.Item(99, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value2 = 2930
$ws.Cells.Item(126, 9).Value2 = 2930
$ws.Cells.Item(126, 10).Value2 = 0
$ws.Cells.Item(126, 11).Value2 = 8790
$ws.Cells.Item(126, 12).Value2 = 0
$ws.Cells.Item(126, 13).Value2 = -6320
$ws.Cells.Item(126, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value2 = 2502.7083
$ws.Cells.Item(132, 9).Value2 = 2293.3
$ws.Cells.Item(132, 11).Value2 = 6879.900000000001
$ws.Cells.Item(132, 13).Value2 = -4349.900000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(104, 8).Value2 = 2020.8889
$ws.Cells.Item(104, 10).Value2 = 2020.8889
$ws.Cells.Item(104, 12).Value2 = 6062.6667
$ws.Cells.Item(104, 14).Value2 = -11304.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value2 = 705.2222
$ws.Cells.Item(122, 9).Value2 = 394.44446
$ws.Cells.Item(122, 10).Value2 = 1016
$ws.Cells.Item(122, 11).Value2 = 3550.00014
$ws.Cells.Item(122, 12).Value2 = 9144
$ws.Cells.Item(122, 13).Value2 = -1100.00014
$ws.Cells.Item(122, 14).Value2 = -14044

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value2 = 30939.281
$ws.Cells.Item(131, 9).Value2 = 1761.3636
$ws.Cells.Item(131, 10).Value2 = 46222.953
$ws.Cells.Item(131, 11).Value2 = 5284.0908
$ws.Cells.Item(131, 12).Value2 = 138668.859
$ws.Cells.Item(131, 13).Value2 = -244.0907999999999
$ws.Cells.Item(131, 14).Value2 = -148748.859

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value2 = 4904.091
$ws.Cells.Item(132, 9).Value2 = 5343.375
$ws.Cells.Item(132, 10).Value2 = 3732.6667
$ws.Cells.Item(132, 11).Value2 = 16030.125
$ws.Cells.Item(132, 12).Value2 = 11198.0001
$ws.Cells.Item(132, 13).Value2 = -13500.125
$ws.Cells.Item(132, 14).Value2 = -16258.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value2 = 20330.75
$ws.Cells.Item(136, 10).Value2 = 20330.75
$ws.Cells.Item(136, 12).Value2 = 60992.25
$ws.Cells.Item(136, 14).Value2 = -66092.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value2 = 4201.4375
$ws.Cells.Item(7, 9).Value2 = 3736.4443
$ws.Cells.Item(7, 10).Value2 = 4799.2856
$ws.Cells.Item(7, 11).Value2 = 3736.4443
$ws.Cells.Item(7, 12).Value2 = 4799.2856
$ws.Cells.Item(7, 13).Value2 = -3624.4443
$ws.Cells.Item(7, 14).Value2 = -5023.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value2 = 769.44446
$ws.Cells.Item(22, 9).Value2 = 500
$ws.Cells.Item(22, 10).Value2 = 803.125
$ws.Cells.Item(22, 11).Value2 = 500
$ws.Cells.Item(22, 12).Value2 = 803.125
$ws.Cells.Item(22, 13).Value2 = -205
$ws.Cells.Item(22, 14).Value2 = -1393.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value2 = 769.44446
$ws.Cells.Item(27, 9).Value2 = 500
$ws.Cells.Item(27, 10).Value2 = 803.125
$ws.Cells.Item(27, 11).Value2 = 500
$ws.Cells.Item(27, 12).Value2 = 803.125
$ws.Cells.Item(27, 13).Value2 = -393
$ws.Cells.Item(27, 14).Value2 = -1017.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value2 = 3733
$ws.Cells.Item(40, 9).Value2 = 3534.08
$ws.Cells.Item(40, 10).Value2 = 4976.25
$ws.Cells.Item(40, 11).Value2 = 3534.08
$ws.Cells.Item(40, 12).Value2 = 4976.25
$ws.Cells.Item(40, 13).Value2 = -3398.08
$ws.Cells.Item(40, 14).Value2 = -5248.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value2 = 799.0909
$ws.Cells.Item(46, 9).Value2 = 648
$ws.Cells.Item(46, 10).Value2 = 855.75
$ws.Cells.Item(46, 11).Value2 = 648
$ws.Cells.Item(46, 12).Value2 = 855.75
$ws.Cells.Item(46, 13).Value2 = -460
$ws.Cells.Item(46, 14).Value2 = -1231.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value2 = 4201.4375
$ws.Cells.Item(126, 9).Value2 = 3736.4443
$ws.Cells.Item(126, 10).Value2 = 4799.2856
$ws.Cells.Item(126, 11).Value2 = 11209.3329
$ws.Cells.Item(126, 12).Value2 = 14397.8568
$ws.Cells.Item(126, 13).Value2 = -8739.332900000001
$ws.Cells.Item(126, 14).Value2 = -19337.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value2 = 6800.086
$ws.Cells.Item(136, 9).Value2 = 5837.6875
$ws.Cells.Item(136, 11).Value2 = 17513.0625
$ws.Cells.Item(136, 13).Value2 = -14963.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value2 = 3840.628
$ws.Cells.Item(136, 9).Value2 = 2922.087
$ws.Cells.Item(136, 10).Value2 = 4896.95
$ws.Cells.Item(136, 11).Value2 = 8766.261
$ws.Cells.Item(136, 12).Value2 = 14690.85
$ws.Cells.Item(136, 13).Value2 = -6216.261
$ws.Cells.Item(136, 14).Value2 = -19790.85
